# Auto-generated Excel COM-interop edit script
# Applies the cryptos.xlsx price/volume/coin updates described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number (e.g. "591.06") need to be
# force-formatted as Text first, otherwise Excel auto-converts the assigned
# string into a numeric value. We flip the format to Text, set the value,
# then restore the cell style to Normal so no numeric style sticks around.
$numericLookingAddrs = @(
    "D5", "D6", "D10", "D12", "D13", "D14", "D16", "D19", "D22", "D23",
    "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33",
    "D35", "D37", "D38", "D39", "D41", "D42", "D44", "D45", "D46", "D47",
    "D48", "D50"
)
foreach ($addr in $numericLookingAddrs) {
    $ws.Range($addr).NumberFormat = "@"
}

# --- Cell value updates ---
$ws.Range("D2").Value = "70.987.76"
$ws.Range("E2").Value = "  +5.66%  "
$ws.Range("D3").Value = "3.638.37"
$ws.Range("E3").Value = "  +5.23%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "591.06"
$ws.Range("E5").Value = "  +2.24%  "
$ws.Range("D6").Value = "194.82"
$ws.Range("E6").Value = "  +3.86%  "
$ws.Range("E7").Value = "  +2.11%  "
$ws.Range("D8").Value = "3.639.75"
$ws.Range("E8").Value = "  +5.71%  "
$ws.Range("E9").Value = "  +0.06%  "
$ws.Range("D10").Value = "0.185"
$ws.Range("E10").Value = "  +7.47%  "
$ws.Range("E11").Value = "  +4.78%  "
$ws.Range("D12").Value = "58.01"
$ws.Range("E12").Value = "  +0.10%  "
$ws.Range("D13").Value = "0.0000312"
$ws.Range("E13").Value = "  +12.53%  "
$ws.Range("D14").Value = "9.91"
$ws.Range("E14").Value = "  +4.50%  "
$ws.Range("D15").Value = "4.229.26"
$ws.Range("E15").Value = "  +5.72%  "
$ws.Range("D16").Value = "20.40"
$ws.Range("E16").Value = "  +7.62%  "
$ws.Range("D17").Value = "3.636.67"
$ws.Range("E17").Value = "  +5.48%  "
$ws.Range("D18").Value = "70.973.60"
$ws.Range("E18").Value = "  +5.79%  "
$ws.Range("D19").Value = "12.77"
$ws.Range("E19").Value = "  +5.88%  "
$ws.Range("E21").Value = "  +3.85%  "
$ws.Range("D22").Value = "491.11"
$ws.Range("E22").Value = "  +0.91%  "
$ws.Range("D23").Value = "18.85"
$ws.Range("E23").Value = "  +9.51%  "
$ws.Range("D24").Value = "5.17"
$ws.Range("E24").Value = "  -5.85%  "
$ws.Range("D25").Value = "4.51"
$ws.Range("E25").Value = "  +3.78%  "
$ws.Range("D26").Value = "91.03"
$ws.Range("E26").Value = "  +1.72%  "
$ws.Range("D27").Value = "3.16"
$ws.Range("E27").Value = "  +6.69%  "
$ws.Range("D28").Value = "11.44"
$ws.Range("E28").Value = "  +4.36%  "
$ws.Range("D29").Value = "9.58"
$ws.Range("E29").Value = "  +6.14%  "
$ws.Range("D30").Value = "7.95"
$ws.Range("E30").Value = "  +7.48%  "
$ws.Range("D31").Value = "32.81"
$ws.Range("E31").Value = "  +4.93%  "
$ws.Range("D32").Value = "0.121"
$ws.Range("E32").Value = "  +8.39%  "
$ws.Range("D33").Value = "67.73"
$ws.Range("E33").Value = "  +4.41%  "
$ws.Range("E34").Value = "  +3.91%  "
$ws.Range("D35").Value = "613.93"
$ws.Range("E35").Value = "  +1.96%  "
$ws.Range("B36").Value = "PEPE"
$ws.Range("C36").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D36").Value = "0.0₃0867"
$ws.Range("E36").Value = "  +9.36%  "
$ws.Range("B37").Value = "InjectiveProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D37").Value = "40.42"
$ws.Range("E37").Value = "  +9.70%  "
$ws.Range("D38").Value = "0.410"
$ws.Range("E38").Value = "  +6.16%  "
$ws.Range("D39").Value = "0.149"
$ws.Range("E39").Value = "  +2.40%  "
$ws.Range("E40").Value = "  -0.02%  "
$ws.Range("B41").Value = "dogwifhat"
$ws.Range("C41").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D41").Value = "3.30"
$ws.Range("E41").Value = "  +23.63%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "3.58"
$ws.Range("E42").Value = "  +3.43%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "3.324.51"
$ws.Range("E43").Value = "  +3.98%  "
$ws.Range("B44").Value = "Fetch.AI"
$ws.Range("C44").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D44").Value = "2.93"
$ws.Range("E44").Value = "  +14.70%  "
$ws.Range("D45").Value = "3.13"
$ws.Range("E45").Value = "  +8.95%  "
$ws.Range("D46").Value = "0.0458"
$ws.Range("E46").Value = "  +6.64%  "
$ws.Range("D47").Value = "9.67"
$ws.Range("E47").Value = "  +11.54%  "
$ws.Range("D48").Value = "3.42"
$ws.Range("E48").Value = "  +6.04%  "
$ws.Range("E49").Value = "  +2.85%  "
$ws.Range("D50").Value = "0.998"
$ws.Range("E50").Value = "  +0.03%  "
$ws.Range("E51").Value = "  +1.32%  "

# Restore Normal style on the text-forced cells so no stray number format remains.
foreach ($addr in $numericLookingAddrs) {
    $ws.Range($addr).Style = "Normal"
}
